$wb = $excel.ActiveWorkbook

# --- Step 1: Insert a new sheet "2022-Q4" right after "总计" ---
$total = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row for the new sheet (style s=2: bold + border + center/top alignment)
$q4.Range("B1:H1").NumberFormat = "@"
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Data rows: column A carries the bold/bordered "index" style (s=2); B:G are
# plain text cells, H is a plain number.
$q4.Range("A2:A27").NumberFormat = "General"
$q4.Range("B2:G27").NumberFormat = "@"
$q4.Range("H2:H27").NumberFormat = "General"

    $q4.Cells.Item(2,1).Value = 0
    $q4.Cells.Item(2,2).Value = "516970"
    $q4.Cells.Item(2,3).Value = "广发中证基建工程ETF"
    $q4.Cells.Item(2,4).Value = "73.53"
    $q4.Cells.Item(2,5).Value = "99.57"
    $q4.Cells.Item(2,6).Value = "3.62"
    $q4.Cells.Item(2,7).Value = "2.6618"
    $q4.Cells.Item(2,8).Value = 9
    $q4.Cells.Item(3,1).Value = 1
    $q4.Cells.Item(3,2).Value = "007202"
    $q4.Cells.Item(3,3).Value = "天弘优质成长企业精选灵活配置混合型证券投资A"
    $q4.Cells.Item(3,4).Value = "5.91"
    $q4.Cells.Item(3,5).Value = "90.29"
    $q4.Cells.Item(3,6).Value = "8.47"
    $q4.Cells.Item(3,7).Value = "0.5006"
    $q4.Cells.Item(3,8).Value = 1
    $q4.Cells.Item(4,1).Value = 2
    $q4.Cells.Item(4,2).Value = "165525"
    $q4.Cells.Item(4,3).Value = "信诚中证基建工程指数（LOF）"
    $q4.Cells.Item(4,4).Value = "9.51"
    $q4.Cells.Item(4,5).Value = "94.24"
    $q4.Cells.Item(4,6).Value = "3.43"
    $q4.Cells.Item(4,7).Value = "0.3262"
    $q4.Cells.Item(4,8).Value = 9
    $q4.Cells.Item(5,1).Value = 3
    $q4.Cells.Item(5,2).Value = "420005"
    $q4.Cells.Item(5,3).Value = "天弘周期策略混合A"
    $q4.Cells.Item(5,4).Value = "3.12"
    $q4.Cells.Item(5,5).Value = "93.31"
    $q4.Cells.Item(5,6).Value = "8.61"
    $q4.Cells.Item(5,7).Value = "0.2686"
    $q4.Cells.Item(5,8).Value = 1
    $q4.Cells.Item(6,1).Value = 4
    $q4.Cells.Item(6,2).Value = "420001"
    $q4.Cells.Item(6,3).Value = "天弘精选混合A"
    $q4.Cells.Item(6,4).Value = "5.36"
    $q4.Cells.Item(6,5).Value = "74.66"
    $q4.Cells.Item(6,6).Value = "4.73"
    $q4.Cells.Item(6,7).Value = "0.2535"
    $q4.Cells.Item(6,8).Value = 3
    $q4.Cells.Item(7,1).Value = 5
    $q4.Cells.Item(7,2).Value = "515150"
    $q4.Cells.Item(7,3).Value = "富国中证国企一带一路ETF"
    $q4.Cells.Item(7,4).Value = "6.96"
    $q4.Cells.Item(7,5).Value = "99.52"
    $q4.Cells.Item(7,6).Value = "2.09"
    $q4.Cells.Item(7,7).Value = "0.1455"
    $q4.Cells.Item(7,8).Value = 8
    $q4.Cells.Item(8,1).Value = 6
    $q4.Cells.Item(8,2).Value = "006022"
    $q4.Cells.Item(8,3).Value = "富国大盘价值量化精选混合A"
    $q4.Cells.Item(8,4).Value = "4.89"
    $q4.Cells.Item(8,5).Value = "93.71"
    $q4.Cells.Item(8,6).Value = "2.75"
    $q4.Cells.Item(8,7).Value = "0.1345"
    $q4.Cells.Item(8,8).Value = 6
    $q4.Cells.Item(9,1).Value = 7
    $q4.Cells.Item(9,2).Value = "011851"
    $q4.Cells.Item(9,3).Value = "天弘先进制造混合A"
    $q4.Cells.Item(9,4).Value = "2.23"
    $q4.Cells.Item(9,5).Value = "84.28"
    $q4.Cells.Item(9,6).Value = "5.80"
    $q4.Cells.Item(9,7).Value = "0.1293"
    $q4.Cells.Item(9,8).Value = 2
    $q4.Cells.Item(10,1).Value = 8
    $q4.Cells.Item(10,2).Value = "001416"
    $q4.Cells.Item(10,3).Value = "嘉实事件驱动股票"
    $q4.Cells.Item(10,4).Value = "9.67"
    $q4.Cells.Item(10,5).Value = "84.39"
    $q4.Cells.Item(10,6).Value = "1.20"
    $q4.Cells.Item(10,7).Value = "0.1160"
    $q4.Cells.Item(10,8).Value = 9
    $q4.Cells.Item(11,1).Value = 9
    $q4.Cells.Item(11,2).Value = "013082"
    $q4.Cells.Item(11,3).Value = "信诚中证基建工程指数（LOF）C"
    $q4.Cells.Item(11,4).Value = "2.86"
    $q4.Cells.Item(11,5).Value = "94.24"
    $q4.Cells.Item(11,6).Value = "3.43"
    $q4.Cells.Item(11,7).Value = "0.0981"
    $q4.Cells.Item(11,8).Value = 9
    $q4.Cells.Item(12,1).Value = 10
    $q4.Cells.Item(12,2).Value = "515110"
    $q4.Cells.Item(12,3).Value = "易方达中证国企一带一路ETF"
    $q4.Cells.Item(12,4).Value = "3.50"
    $q4.Cells.Item(12,5).Value = "98.83"
    $q4.Cells.Item(12,6).Value = "2.09"
    $q4.Cells.Item(12,7).Value = "0.0732"
    $q4.Cells.Item(12,8).Value = 8
    $q4.Cells.Item(13,1).Value = 11
    $q4.Cells.Item(13,2).Value = "015458"
    $q4.Cells.Item(13,3).Value = "天弘周期策略混合C"
    $q4.Cells.Item(13,4).Value = "0.56"
    $q4.Cells.Item(13,5).Value = "93.31"
    $q4.Cells.Item(13,6).Value = "8.61"
    $q4.Cells.Item(13,7).Value = "0.0482"
    $q4.Cells.Item(13,8).Value = 1
    $q4.Cells.Item(14,1).Value = 12
    $q4.Cells.Item(14,2).Value = "011852"
    $q4.Cells.Item(14,3).Value = "天弘先进制造混合C"
    $q4.Cells.Item(14,4).Value = "0.67"
    $q4.Cells.Item(14,5).Value = "84.28"
    $q4.Cells.Item(14,6).Value = "5.80"
    $q4.Cells.Item(14,7).Value = "0.0389"
    $q4.Cells.Item(14,8).Value = 2
    $q4.Cells.Item(15,1).Value = 13
    $q4.Cells.Item(15,2).Value = "004694"
    $q4.Cells.Item(15,3).Value = "天弘策略精选灵活配置混合A"
    $q4.Cells.Item(15,4).Value = "0.82"
    $q4.Cells.Item(15,5).Value = "87.95"
    $q4.Cells.Item(15,6).Value = "3.69"
    $q4.Cells.Item(15,7).Value = "0.0303"
    $q4.Cells.Item(15,8).Value = 4
    $q4.Cells.Item(16,1).Value = 14
    $q4.Cells.Item(16,2).Value = "012879"
    $q4.Cells.Item(16,3).Value = "中信建投量化精选6个月持有期混合C"
    $q4.Cells.Item(16,4).Value = "3.12"
    $q4.Cells.Item(16,5).Value = "70.35"
    $q4.Cells.Item(16,6).Value = "0.84"
    $q4.Cells.Item(16,7).Value = "0.0262"
    $q4.Cells.Item(16,8).Value = 6
    $q4.Cells.Item(17,1).Value = 15
    $q4.Cells.Item(17,2).Value = "515990"
    $q4.Cells.Item(17,3).Value = "汇添富中证国企一带一路ETF"
    $q4.Cells.Item(17,4).Value = "0.97"
    $q4.Cells.Item(17,5).Value = "98.95"
    $q4.Cells.Item(17,6).Value = "2.09"
    $q4.Cells.Item(17,7).Value = "0.0203"
    $q4.Cells.Item(17,8).Value = 8
    $q4.Cells.Item(18,1).Value = 16
    $q4.Cells.Item(18,2).Value = "009188"
    $q4.Cells.Item(18,3).Value = "鹏华股息精选混合"
    $q4.Cells.Item(18,4).Value = "0.72"
    $q4.Cells.Item(18,5).Value = "87.00"
    $q4.Cells.Item(18,6).Value = "1.86"
    $q4.Cells.Item(18,7).Value = "0.0134"
    $q4.Cells.Item(18,8).Value = 8
    $q4.Cells.Item(19,1).Value = 17
    $q4.Cells.Item(19,2).Value = "012878"
    $q4.Cells.Item(19,3).Value = "中信建投量化精选6个月持有期混合A"
    $q4.Cells.Item(19,4).Value = "1.59"
    $q4.Cells.Item(19,5).Value = "70.35"
    $q4.Cells.Item(19,6).Value = "0.84"
    $q4.Cells.Item(19,7).Value = "0.0134"
    $q4.Cells.Item(19,8).Value = 6
    $q4.Cells.Item(20,1).Value = 18
    $q4.Cells.Item(20,2).Value = "015460"
    $q4.Cells.Item(20,3).Value = "天弘优质成长企业精选灵活配置混合型证券投资C"
    $q4.Cells.Item(20,4).Value = "0.13"
    $q4.Cells.Item(20,5).Value = "90.29"
    $q4.Cells.Item(20,6).Value = "8.47"
    $q4.Cells.Item(20,7).Value = "0.0110"
    $q4.Cells.Item(20,8).Value = 1
    $q4.Cells.Item(21,1).Value = 19
    $q4.Cells.Item(21,2).Value = "410009"
    $q4.Cells.Item(21,3).Value = "华富量子生命力混合"
    $q4.Cells.Item(21,4).Value = "0.10"
    $q4.Cells.Item(21,5).Value = "92.12"
    $q4.Cells.Item(21,6).Value = "5.08"
    $q4.Cells.Item(21,7).Value = "0.0051"
    $q4.Cells.Item(21,8).Value = 4
    $q4.Cells.Item(22,1).Value = 20
    $q4.Cells.Item(22,2).Value = "009384"
    $q4.Cells.Item(22,3).Value = "摩根士丹利华鑫MSCI中国A股指数增强A"
    $q4.Cells.Item(22,4).Value = "0.38"
    $q4.Cells.Item(22,5).Value = "89.37"
    $q4.Cells.Item(22,6).Value = "1.11"
    $q4.Cells.Item(22,7).Value = "0.0042"
    $q4.Cells.Item(22,8).Value = 6
    $q4.Cells.Item(23,1).Value = 21
    $q4.Cells.Item(23,2).Value = "004748"
    $q4.Cells.Item(23,3).Value = "天弘策略精选灵活配置混合C"
    $q4.Cells.Item(23,4).Value = "0.05"
    $q4.Cells.Item(23,5).Value = "87.95"
    $q4.Cells.Item(23,6).Value = "3.69"
    $q4.Cells.Item(23,7).Value = "0.0018"
    $q4.Cells.Item(23,8).Value = 4
    $q4.Cells.Item(24,1).Value = 22
    $q4.Cells.Item(24,2).Value = "007808"
    $q4.Cells.Item(24,3).Value = "北信瑞丰量化优选灵活配置混合"
    $q4.Cells.Item(24,4).Value = "0.15"
    $q4.Cells.Item(24,5).Value = "79.42"
    $q4.Cells.Item(24,6).Value = "1.15"
    $q4.Cells.Item(24,7).Value = "0.0017"
    $q4.Cells.Item(24,8).Value = 6
    $q4.Cells.Item(25,1).Value = 23
    $q4.Cells.Item(25,2).Value = "014866"
    $q4.Cells.Item(25,3).Value = "摩根士丹利华鑫MSCI中国A股指数增强C"
    $q4.Cells.Item(25,4).Value = "0.00"
    $q4.Cells.Item(25,5).Value = "89.37"
    $q4.Cells.Item(25,6).Value = "1.11"
    $q4.Cells.Item(25,7).NumberFormat = "General"
    $q4.Cells.Item(25,7).Value = 0
    $q4.Cells.Item(25,8).Value = 6
    $q4.Cells.Item(26,1).Value = 24
    $q4.Cells.Item(26,2).Value = "015459"
    $q4.Cells.Item(26,3).Value = "天弘精选混合C"
    $q4.Cells.Item(26,4).Value = "0.00"
    $q4.Cells.Item(26,5).Value = "74.66"
    $q4.Cells.Item(26,6).Value = "4.73"
    $q4.Cells.Item(26,7).NumberFormat = "General"
    $q4.Cells.Item(26,7).Value = 0
    $q4.Cells.Item(26,8).Value = 3
    $q4.Cells.Item(27,1).Value = 25
    $q4.Cells.Item(27,2).Value = "014181"
    $q4.Cells.Item(27,3).Value = "富国大盘价值量化精选混合C"
    $q4.Cells.Item(27,4).Value = "0.00"
    $q4.Cells.Item(27,5).Value = "93.71"
    $q4.Cells.Item(27,6).Value = "2.75"
    $q4.Cells.Item(27,7).NumberFormat = "General"
    $q4.Cells.Item(27,7).Value = 0
    $q4.Cells.Item(27,8).Value = 6

# --- Step 2: copy the header/index styling (s=2) from an existing sheet ---
# Every quarterly sheet shares the same bold+bordered style for its header
# row (B1:H1) and its index column (A2:A.. ); grab it from the "2022-Q3"
# sheet (now shifted down to make room for the new one) and tile it onto the
# freshly-created "2022-Q4" sheet.
$src = $wb.Worksheets.Item("2022-Q3")
$src.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2").Copy()
$q4.Range("A2:A27").PasteSpecial(-4122)
$excel.CutCopyMode = $false


# --- Step 3: update the "总计" (summary) sheet ---
# Insert a new row 2 for "2022-Q4" and push the existing quarters down.
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 26
$total.Cells.Item(2,4).Value = 4.92

# Restore the bold+bordered style on the new index cell (A2) and renumber
# the index column for every row pushed down (it is a simple 0-based
# position counter, not data that shifts verbatim).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
for ($r = 3; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
